$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 12

# Copy the style of the previous year's label cell (bold / bordered / centered)
# onto the new label cell before setting its value.
$ws.Cells.Item($row - 1, 1).Copy()
$ws.Cells.Item($row, 1).PasteSpecial(-4122)

$ws.Cells.Item($row, 1).Value = "2021年"

$ws.Cells.Item($row, 2).Value = 5234.9
$ws.Cells.Item($row, 3).Value = 8844.4
$ws.Cells.Item($row, 4).Value = 216238.6
$ws.Cells.Item($row, 5).Value = ""
$ws.Cells.Item($row, 6).Value = 1049041.2
$ws.Cells.Item($row, 7).Value = 1579449.6
$ws.Cells.Item($row, 8).Value = 432207.1
$ws.Cells.Item($row, 9).Value = 248485.3
$ws.Cells.Item($row, 10).Value = 248292.8
$ws.Cells.Item($row, 11).Value = 600048
$ws.Cells.Item($row, 12).Value = 1377176.7
$ws.Cells.Item($row, 13).Value = ""
$ws.Cells.Item($row, 14).Value = 33072
$ws.Cells.Item($row, 15).Value = 144054.9
$ws.Cells.Item($row, 16).Value = 352928.9
$ws.Cells.Item($row, 17).Value = 1977224.7
$ws.Cells.Item($row, 18).Value = 3809.5
$ws.Cells.Item($row, 19).Value = 68806.8
